# Clear/update computed result columns on Sheet1 (rows 1-16, columns A:AD)
# per updated model parameters (see commit: "modify parameter for clear result").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object "object[,]" 16,30

# Row 1
$arr[0,0] = [double]"1"
$arr[0,1] = [double]"0"
$arr[0,2] = [double]"0"
$arr[0,3] = [double]"0"
$arr[0,4] = [double]"0"
$arr[0,5] = [double]"0"
$arr[0,6] = [double]"0"
$arr[0,7] = [double]"0"
$arr[0,8] = [double]"0"
$arr[0,9] = [double]"76.570724343257581"
$arr[0,10] = [double]"0"
$arr[0,11] = [double]"0"
$arr[0,12] = [double]"6.3064123963058666E-3"
$arr[0,13] = [double]"49.471903250931668"
$arr[0,14] = [double]"0"
$arr[0,15] = [double]"0"
$arr[0,16] = [double]"76.570724343257581"
$arr[0,17] = [double]"0"
$arr[0,18] = [double]"0"
$arr[0,19] = [double]"1.0059653637984452E-2"
$arr[0,20] = [double]"49.471903250931668"
$arr[0,21] = [double]"0"
$arr[0,22] = [double]"0"
$arr[0,23] = [double]"76.570724343257581"
$arr[0,24] = [double]"0"
$arr[0,25] = [double]"0"
$arr[0,26] = [double]"3.1286871079647789E-2"
$arr[0,27] = [double]"49.471903250931668"
$arr[0,28] = [double]"0"
$arr[0,29] = [double]"0"

# Row 2
$arr[1,0] = [double]"2"
$arr[1,1] = [double]"0"
$arr[1,2] = [double]"0"
$arr[1,3] = [double]"0"
$arr[1,4] = [double]"0"
$arr[1,5] = [double]"0"
$arr[1,6] = [double]"0"
$arr[1,7] = [double]"0"
$arr[1,8] = [double]"0"
$arr[1,9] = [double]"2032.5822270573306"
$arr[1,10] = [double]"0"
$arr[1,11] = [double]"0"
$arr[1,12] = [double]"1.6018759663991754E-2"
$arr[1,13] = [double]"2007.005405965006"
$arr[1,14] = [double]"0"
$arr[1,15] = [double]"0"
$arr[1,16] = [double]"209.81118823065617"
$arr[1,17] = [double]"0"
$arr[1,18] = [double]"0"
$arr[1,19] = [double]"1.8372819877415927E-2"
$arr[1,20] = [double]"121.8368094235322"
$arr[1,21] = [double]"0"
$arr[1,22] = [double]"0"
$arr[1,23] = [double]"209.81118823065617"
$arr[1,24] = [double]"0"
$arr[1,25] = [double]"0"
$arr[1,26] = [double]"6.0741068167621828E-2"
$arr[1,27] = [double]"121.8368094235322"
$arr[1,28] = [double]"0"
$arr[1,29] = [double]"0"

# Row 3
$arr[2,0] = [double]"3"
$arr[2,1] = [double]"0"
$arr[2,2] = [double]"0"
$arr[2,3] = [double]"0"
$arr[2,4] = [double]"0"
$arr[2,5] = [double]"0"
$arr[2,6] = [double]"0"
$arr[2,7] = [double]"0"
$arr[2,8] = [double]"0"
$arr[2,9] = [double]"20373.599423022504"
$arr[2,10] = [double]"18182.944399996137"
$arr[2,11] = [double]"2"
$arr[2,12] = [double]"1.4667934150029526E-2"
$arr[2,13] = [double]"20300.953409663798"
$arr[2,14] = [double]"2"
$arr[2,15] = [double]"0"
$arr[2,16] = [double]"18510.346990778009"
$arr[2,17] = [double]"18182.944399996137"
$arr[2,18] = [double]"2"
$arr[2,19] = [double]"1.71738875671715E-2"
$arr[2,20] = [double]"18419.20281312232"
$arr[2,21] = [double]"2"
$arr[2,22] = [double]"0"
$arr[2,23] = [double]"892.29741555146859"
$arr[2,24] = [double]"0"
$arr[2,25] = [double]"0"
$arr[2,26] = [double]"2.3998557299008884E-2"
$arr[2,27] = [double]"418.32270975413138"
$arr[2,28] = [double]"0"
$arr[2,29] = [double]"0"

# Row 4
$arr[3,0] = [double]"4"
$arr[3,1] = [double]"0"
$arr[3,2] = [double]"0"
$arr[3,3] = [double]"0"
$arr[3,4] = [double]"0"
$arr[3,5] = [double]"0"
$arr[3,6] = [double]"0"
$arr[3,7] = [double]"0"
$arr[3,8] = [double]"0"
$arr[3,9] = [double]"20915.296883265106"
$arr[3,10] = [double]"18536.466660293354"
$arr[3,11] = [double]"2"
$arr[3,12] = [double]"3.5543863002570238E-2"
$arr[3,13] = [double]"20802.088707995001"
$arr[3,14] = [double]"2"
$arr[3,15] = [double]"0"
$arr[3,16] = [double]"19244.538922948941"
$arr[3,17] = [double]"18536.466660293354"
$arr[3,18] = [double]"2"
$arr[3,19] = [double]"2.3055624770311308E-2"
$arr[3,20] = [double]"18921.309111453545"
$arr[3,21] = [double]"2"
$arr[3,22] = [double]"0"
$arr[3,23] = [double]"532.63102317939195"
$arr[3,24] = [double]"0"
$arr[3,25] = [double]"0"
$arr[3,26] = [double]"4.1771768799201733E-2"
$arr[3,27] = [double]"390.33571677624673"
$arr[3,28] = [double]"0"
$arr[3,29] = [double]"0"

# Row 5
$arr[4,0] = [double]"5"
$arr[4,1] = [double]"0"
$arr[4,2] = [double]"0"
$arr[4,3] = [double]"0"
$arr[4,4] = [double]"0"
$arr[4,5] = [double]"0"
$arr[4,6] = [double]"0"
$arr[4,7] = [double]"0"
$arr[4,8] = [double]"0"
$arr[4,9] = [double]"21408.914923198958"
$arr[4,10] = [double]"18937.125221963535"
$arr[4,11] = [double]"2"
$arr[4,12] = [double]"4.4440993188130259E-2"
$arr[4,13] = [double]"21156.25512321986"
$arr[4,14] = [double]"2"
$arr[4,15] = [double]"0"
$arr[4,16] = [double]"19853.217784146193"
$arr[4,17] = [double]"18937.125221963535"
$arr[4,18] = [double]"2"
$arr[4,19] = [double]"2.8037664963414786E-2"
$arr[4,20] = [double]"19280.140526678369"
$arr[4,21] = [double]"2"
$arr[4,22] = [double]"0"
$arr[4,23] = [double]"630.70438594655207"
$arr[4,24] = [double]"0"
$arr[4,25] = [double]"0"
$arr[4,26] = [double]"5.2773786077480324E-2"
$arr[4,27] = [double]"466.14830888783263"
$arr[4,28] = [double]"0"
$arr[4,29] = [double]"0"

# Row 6
$arr[5,0] = [double]"6"
$arr[5,1] = [double]"0"
$arr[5,2] = [double]"0"
$arr[5,3] = [double]"0"
$arr[5,4] = [double]"0"
$arr[5,5] = [double]"0"
$arr[5,6] = [double]"0"
$arr[5,7] = [double]"0"
$arr[5,8] = [double]"0"
$arr[5,9] = [double]"21945.177342471412"
$arr[5,10] = [double]"19395.020721015175"
$arr[5,11] = [double]"2"
$arr[5,12] = [double]"4.3475020679093465E-2"
$arr[5,13] = [double]"21676.381264704567"
$arr[5,14] = [double]"2"
$arr[5,15] = [double]"0"
$arr[5,16] = [double]"20437.49430504446"
$arr[5,17] = [double]"19395.020721015175"
$arr[5,18] = [double]"2"
$arr[5,19] = [double]"5.1096418619945003E-2"
$arr[5,20] = [double]"19794.551668163065"
$arr[5,21] = [double]"2"
$arr[5,22] = [double]"0"
$arr[5,23] = [double]"871.70748500309378"
$arr[5,24] = [double]"0"
$arr[5,25] = [double]"0"
$arr[5,26] = [double]"6.1822525022556425E-2"
$arr[5,27] = [double]"730.25254414923518"
$arr[5,28] = [double]"0"
$arr[5,29] = [double]"0"

# Row 7
$arr[6,0] = [double]"7"
$arr[6,1] = [double]"0"
$arr[6,2] = [double]"0"
$arr[6,3] = [double]"0"
$arr[6,4] = [double]"0"
$arr[6,5] = [double]"0"
$arr[6,6] = [double]"0"
$arr[6,7] = [double]"0"
$arr[6,8] = [double]"0"
$arr[6,9] = [double]"55328.337110708431"
$arr[6,10] = [double]"52685.207835165646"
$arr[6,11] = [double]"4"
$arr[6,12] = [double]"3.6970350674189648E-2"
$arr[6,13] = [double]"55045.102794036102"
$arr[6,14] = [double]"4"
$arr[6,15] = [double]"0"
$arr[6,16] = [double]"53833.42940173544"
$arr[6,17] = [double]"52685.207835165653"
$arr[6,18] = [double]"4"
$arr[6,19] = [double]"3.9443308563932244E-2"
$arr[6,20] = [double]"53160.90119749459"
$arr[6,21] = [double]"4"
$arr[6,22] = [double]"0"
$arr[6,23] = [double]"1014.1532551878734"
$arr[6,24] = [double]"0"
$arr[6,25] = [double]"0"
$arr[6,26] = [double]"0.32835356869384358"
$arr[6,27] = [double]"782.9718564620332"
$arr[6,28] = [double]"0"
$arr[6,29] = [double]"0"

# Row 8
$arr[7,0] = [double]"8"
$arr[7,1] = [double]"0"
$arr[7,2] = [double]"0"
$arr[7,3] = [double]"0"
$arr[7,4] = [double]"0"
$arr[7,5] = [double]"0"
$arr[7,6] = [double]"0"
$arr[7,7] = [double]"0"
$arr[7,8] = [double]"0"
$arr[7,9] = [double]"60567.736133916282"
$arr[7,10] = [double]"57781.636584276021"
$arr[7,11] = [double]"5"
$arr[7,12] = [double]"4.5575925552987971E-2"
$arr[7,13] = [double]"60226.4349094055"
$arr[7,14] = [double]"5"
$arr[7,15] = [double]"0"
$arr[7,16] = [double]"58994.062532449818"
$arr[7,17] = [double]"57781.636584276021"
$arr[7,18] = [double]"5"
$arr[7,19] = [double]"4.7037399861418784E-2"
$arr[7,20] = [double]"58345.031312864012"
$arr[7,21] = [double]"5"
$arr[7,22] = [double]"0"
$arr[7,23] = [double]"1095.7786130443003"
$arr[7,24] = [double]"0"
$arr[7,25] = [double]"0"
$arr[7,26] = [double]"0.23251265492253448"
$arr[7,27] = [double]"814.33781538367464"
$arr[7,28] = [double]"0"
$arr[7,29] = [double]"0"

# Row 9
$arr[8,0] = [double]"9"
$arr[8,1] = [double]"0"
$arr[8,2] = [double]"0"
$arr[8,3] = [double]"0"
$arr[8,4] = [double]"0"
$arr[8,5] = [double]"0"
$arr[8,6] = [double]"0"
$arr[8,7] = [double]"0"
$arr[8,8] = [double]"0"
$arr[8,9] = [double]"62483.26710327663"
$arr[8,10] = [double]"59618.170391417334"
$arr[8,11] = [double]"5"
$arr[8,12] = [double]"5.7682724110742094E-2"
$arr[8,13] = [double]"62118.431921620882"
$arr[8,14] = [double]"5"
$arr[8,15] = [double]"0"
$arr[8,16] = [double]"60943.880277231852"
$arr[8,17] = [double]"59618.170391417334"
$arr[8,18] = [double]"5"
$arr[8,19] = [double]"5.0879387694033613E-2"
$arr[8,20] = [double]"60240.592325079422"
$arr[8,21] = [double]"5"
$arr[8,22] = [double]"0"
$arr[8,23] = [double]"1227.4930050061523"
$arr[8,24] = [double]"0"
$arr[8,25] = [double]"0"
$arr[8,26] = [double]"1.0649252423749043"
$arr[8,27] = [double]"874.93968417242274"
$arr[8,28] = [double]"0"
$arr[8,29] = [double]"0"

# Row 10
$arr[9,0] = [double]"10"
$arr[9,1] = [double]"0"
$arr[9,2] = [double]"0"
$arr[9,3] = [double]"0"
$arr[9,4] = [double]"0"
$arr[9,5] = [double]"0"
$arr[9,6] = [double]"0"
$arr[9,7] = [double]"0"
$arr[9,8] = [double]"0"
$arr[9,9] = [double]"71594.86295082829"
$arr[9,10] = [double]"66427.181145974944"
$arr[9,11] = [double]"6"
$arr[9,12] = [double]"5.2305875365653023E-2"
$arr[9,13] = [double]"71227.895769172654"
$arr[9,14] = [double]"6"
$arr[9,15] = [double]"0"
$arr[9,16] = [double]"63263.301030696006"
$arr[9,17] = [double]"61822.010959986859"
$arr[9,18] = [double]"5"
$arr[9,19] = [double]"6.6884209592141142E-2"
$arr[9,20] = [double]"62547.437903997357"
$arr[9,21] = [double]"5"
$arr[9,22] = [double]"0"
$arr[9,23] = [double]"2826.1500983053807"
$arr[9,24] = [double]"0"
$arr[9,25] = [double]"0"
$arr[9,26] = [double]"1.9646083146453492"
$arr[9,27] = [double]"1315.36879153985"
$arr[9,28] = [double]"0"
$arr[9,29] = [double]"0"

# Row 11
$arr[10,0] = [double]"11"
$arr[10,1] = [double]"0"
$arr[10,2] = [double]"0"
$arr[10,3] = [double]"0"
$arr[10,4] = [double]"0"
$arr[10,5] = [double]"0"
$arr[10,6] = [double]"0"
$arr[10,7] = [double]"0"
$arr[10,8] = [double]"0"
$arr[10,9] = [double]"82219.572299392763"
$arr[10,10] = [double]"74749.305401545382"
$arr[10,11] = [double]"7"
$arr[10,12] = [double]"5.4344202515169408E-2"
$arr[10,13] = [double]"81858.440117736987"
$arr[10,14] = [double]"7"
$arr[10,15] = [double]"0"
$arr[10,16] = [double]"66474.616538476053"
$arr[10,17] = [double]"64515.593877127409"
$arr[10,18] = [double]"5"
$arr[10,19] = [double]"6.3872514384343285E-2"
$arr[10,20] = [double]"65349.934634289442"
$arr[10,21] = [double]"5"
$arr[10,22] = [double]"0"
$arr[10,23] = [double]"1592.3628330484526"
$arr[10,24] = [double]"0"
$arr[10,25] = [double]"0"
$arr[10,26] = [double]"7.4793399642965639"
$arr[10,27] = [double]"1345.2175270453586"
$arr[10,28] = [double]"0"
$arr[10,29] = [double]"0"

# Row 12
$arr[11,0] = [double]"12"
$arr[11,1] = [double]"0"
$arr[11,2] = [double]"0"
$arr[11,3] = [double]"0"
$arr[11,4] = [double]"0"
$arr[11,5] = [double]"0"
$arr[11,6] = [double]"0"
$arr[11,7] = [double]"0"
$arr[11,8] = [double]"0"
$arr[11,9] = [double]"88821.7763890848"
$arr[11,10] = [double]"79395.497988523333"
$arr[11,11] = [double]"7"
$arr[11,12] = [double]"5.8377905739742657E-2"
$arr[11,13] = [double]"88445.919207429004"
$arr[11,14] = [double]"7"
$arr[11,15] = [double]"0"
$arr[11,16] = [double]"88735.309131999209"
$arr[11,17] = [double]"86672.454876414064"
$arr[11,18] = [double]"6"
$arr[11,19] = [double]"6.9726376055714648E-2"
$arr[11,20] = [double]"87471.561876539388"
$arr[11,21] = [double]"6"
$arr[11,22] = [double]"0"
$arr[11,23] = [double]"1870.3338046630006"
$arr[11,24] = [double]"0"
$arr[11,25] = [double]"0"
$arr[11,26] = [double]"2.3509321805378725"
$arr[11,27] = [double]"1458.9098576314186"
$arr[11,28] = [double]"0"
$arr[11,29] = [double]"0"

# Row 13
$arr[12,0] = [double]"13"
$arr[12,1] = [double]"0"
$arr[12,2] = [double]"0"
$arr[12,3] = [double]"0"
$arr[12,4] = [double]"0"
$arr[12,5] = [double]"0"
$arr[12,6] = [double]"0"
$arr[12,7] = [double]"0"
$arr[12,8] = [double]"0"
$arr[12,9] = [double]"94919.778452457758"
$arr[12,10] = [double]"85369.174171780687"
$arr[12,11] = [double]"7"
$arr[12,12] = [double]"6.6828458528604279E-2"
$arr[12,13] = [double]"94536.824523701041"
$arr[12,14] = [double]"7"
$arr[12,15] = [double]"0"
$arr[12,16] = [double]"94142.852696075628"
$arr[12,17] = [double]"91896.183724335715"
$arr[12,18] = [double]"6"
$arr[12,19] = [double]"7.4912647185652173E-2"
$arr[12,20] = [double]"92808.291857475488"
$arr[12,21] = [double]"6"
$arr[12,22] = [double]"0"
$arr[12,23] = [double]"94142.852696075628"
$arr[12,24] = [double]"91896.183724335715"
$arr[12,25] = [double]"6"
$arr[12,26] = [double]"9.3393843148364954"
$arr[12,27] = [double]"92808.291857475488"
$arr[12,28] = [double]"6"
$arr[12,29] = [double]"0"

# Row 14
$arr[13,0] = [double]"14"
$arr[13,1] = [double]"0"
$arr[13,2] = [double]"0"
$arr[13,3] = [double]"0"
$arr[13,4] = [double]"0"
$arr[13,5] = [double]"0"
$arr[13,6] = [double]"0"
$arr[13,7] = [double]"0"
$arr[13,8] = [double]"0"
$arr[13,9] = [double]"104840.6915328483"
$arr[13,10] = [double]"93334.075749457159"
$arr[13,11] = [double]"7"
$arr[13,12] = [double]"6.8222235117026034E-2"
$arr[13,13] = [double]"104462.76260409155"
$arr[13,14] = [double]"7"
$arr[13,15] = [double]"0"
$arr[13,16] = [double]"101224.33909343324"
$arr[13,17] = [double]"98861.155521564564"
$arr[13,18] = [double]"6"
$arr[13,19] = [double]"9.8176769555823315E-2"
$arr[13,20] = [double]"99858.264475064701"
$arr[13,21] = [double]"6"
$arr[13,22] = [double]"0"
$arr[13,23] = [double]"101224.33909343324"
$arr[13,24] = [double]"98861.155521564564"
$arr[13,25] = [double]"6"
$arr[13,26] = [double]"5.5232811890222742"
$arr[13,27] = [double]"99858.264475064701"
$arr[13,28] = [double]"6"
$arr[13,29] = [double]"0"

# Row 15
$arr[14,0] = [double]"15"
$arr[14,1] = [double]"0"
$arr[14,2] = [double]"0"
$arr[14,3] = [double]"0"
$arr[14,4] = [double]"0"
$arr[14,5] = [double]"0"
$arr[14,6] = [double]"0"
$arr[14,7] = [double]"0"
$arr[14,8] = [double]"0"
$arr[14,9] = [double]"147839.69232343382"
$arr[14,10] = [double]"136243.1371536135"
$arr[14,11] = [double]"9"
$arr[14,12] = [double]"7.4846371686651716E-2"
$arr[14,13] = [double]"147427.52406154003"
$arr[14,14] = [double]"9"
$arr[14,15] = [double]"0"
$arr[14,16] = [double]"142927.92483813508"
$arr[14,17] = [double]"140370.31523309427"
$arr[14,18] = [double]"8"
$arr[14,19] = [double]"8.7120138990814722E-2"
$arr[14,20] = [double]"141431.0912398865"
$arr[14,21] = [double]"8"
$arr[14,22] = [double]"0"
$arr[14,23] = [double]"142927.92483813508"
$arr[14,24] = [double]"140370.31523309427"
$arr[14,25] = [double]"8"
$arr[14,26] = [double]"3.3074847509463456"
$arr[14,27] = [double]"141431.0912398865"
$arr[14,28] = [double]"8"
$arr[14,29] = [double]"0"

# Row 16
$arr[15,0] = [double]"16"
$arr[15,1] = [double]"0"
$arr[15,2] = [double]"0"
$arr[15,3] = [double]"0"
$arr[15,4] = [double]"0"
$arr[15,5] = [double]"0"
$arr[15,6] = [double]"0"
$arr[15,7] = [double]"0"
$arr[15,8] = [double]"0"
$arr[15,9] = [double]"195258.01963256596"
$arr[15,10] = [double]"183498.59438313707"
$arr[15,11] = [double]"16"
$arr[15,12] = [double]"7.6387774815765488E-2"
$arr[15,13] = [double]"194738.77103953919"
$arr[15,14] = [double]"16"
$arr[15,15] = [double]"0"
$arr[15,16] = [double]"187867.03657092812"
$arr[15,17] = [double]"185191.87757073811"
$arr[15,18] = [double]"12"
$arr[15,19] = [double]"9.1602410721498412E-2"
$arr[15,20] = [double]"186309.6203260062"
$arr[15,21] = [double]"12"
$arr[15,22] = [double]"0"
$arr[15,23] = [double]"187867.03657092812"
$arr[15,24] = [double]"185191.87757073811"
$arr[15,25] = [double]"12"
$arr[15,26] = [double]"1.5924986944011141"
$arr[15,27] = [double]"186309.6203260062"
$arr[15,28] = [double]"12"
$arr[15,29] = [double]"0"

$ws.Range("A1:AD16").Value2 = $arr
Write-Host "Updated rows 1-16 (A:AD) with refreshed model results"
